$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new transaction was recorded on 2024-09-05 16:28:38 ("transfer"),
# which pushes all existing September entries (rows 31-69), the August
# entries (rows 70-73) and the "Broadband" row (74) down by one row.
$ws.Rows.Item(31).Insert()

$ws.Range("R31").Value = "transfer"
$ws.Range("S31").Value = "2024-09-05 16:28:38"
